$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 380; this shifts existing rows 380..397 down
# to 381..398 and copies formatting (e.g. the date style on column D) from
# the row above, matching how the source workbook grew by one weekly
# observation.
$ws.Rows.Item(380).Insert()

# Populate the newly inserted row 380 with the new weekly price observation.
# Columns A, B, C, E, F, G, H, I, N, Q, R are constant for every row in this
# market/product block, so we simply repeat them.
$ws.Cells.Item(380, 1).Value = 3
$ws.Cells.Item(380, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(380, 3).Value = "Coquimbo"
$ws.Cells.Item(380, 4).Value = 44753
$ws.Cells.Item(380, 5).Value = 5
$ws.Cells.Item(380, 6).Value = 100112031
$ws.Cells.Item(380, 7).Value = "Poroto verde"
$ws.Cells.Item(380, 8).Value = "Magnum"
$ws.Cells.Item(380, 9).Value = "Primera"
$ws.Cells.Item(380, 10).Value = 65
$ws.Cells.Item(380, 11).Value = 32000
$ws.Cells.Item(380, 12).Value = 33000
$ws.Cells.Item(380, 13).Value = 32538
$ws.Cells.Item(380, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(380, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(380, 16).Value = 1302
$ws.Cells.Item(380, 17).Value = 25
$ws.Cells.Item(380, 18).Value = "Hortaliza"
